$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update consultation_area for test1 (row 2): "In City" -> "All Over The World"
$ws.Range("D2").Value = "All Over The World"

# Update consultation_type for test2 (row 3): "Both" -> "Online"
$ws.Range("C3").Value = "Online"

# Update the saved selection to match the authored state
$ws.Range("G7").Select()
